# Rename two of the header labels on Sheet1:
#   H1: "six year graduation count"  -> "totalCohortToGradCount"
#   F1: "Graduation Rate"            -> "gradRate"
# (H1 is set first so the new shared-string entries land in the same
#  order as the target workbook: totalCohortToGradCount before gradRate.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "totalCohortToGradCount"
$ws.Range("F1").Value = "gradRate"

# Move the active selection to F1 (was F8).
$ws.Range("F1").Select() | Out-Null
